# approach5 corrected and new summary generated
#
# 1) On "approach5.xlsx" sheet, a new row for "template15" was inserted
#    right after "template14" (row 16), pushing the former rows 16-25
#    (template16..template25) down to rows 17-26 unchanged.
# 2) On "general_summary" sheet, row 6 (the approach5.xlsx summary row)
#    was recomputed with the corrected numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the missing "template15" row into approach5.xlsx
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("approach5.xlsx")

$ws5.Rows.Item(16).Insert()

$ws5.Range("A16").Value = "template15"
$ws5.Range("B16").Value = 14
$ws5.Range("C16").Value = 9
$ws5.Range("D16").Value = 5
$ws5.Range("E16").Value = 0
$ws5.Range("F16").Value = 0
$ws5.Range("G16").Value = 9
$ws5.Range("H16").Value = 5
$ws5.Range("I16").Value = 0
$ws5.Range("J16").Value = 0
$ws5.Range("K16").Value = 0
$ws5.Range("L16").Value = 0
$ws5.Range("M16").Value = 0
$ws5.Range("N16").Value = 0

# ---------------------------------------------------------------------
# Step 2: update the approach5.xlsx summary row on general_summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("general_summary")

$wsSummary.Range("B6").Value = 386
$wsSummary.Range("C6").Value = 230
$wsSummary.Range("E6").Value = 156
$wsSummary.Range("K6").Value = 211
$wsSummary.Range("M6").Value = 93

# columns D, F, H, J, L, N hold percentages stored as text in this
# workbook (e.g. "59.6"), so force text formatting before writing,
# otherwise Excel will auto-convert the numeric-looking text to a number.
# (NumberFormat is applied cell-by-cell: a multi-area Range only applies
# formatting to its first area.)
$wsSummary.Range("D6").NumberFormat = "@"
$wsSummary.Range("D6").Value = "59.6"
$wsSummary.Range("D6").Style = "Normal"

$wsSummary.Range("F6").NumberFormat = "@"
$wsSummary.Range("F6").Value = "40.4"
$wsSummary.Range("F6").Style = "Normal"

$wsSummary.Range("H6").NumberFormat = "@"
$wsSummary.Range("H6").Value = "4.9"
$wsSummary.Range("H6").Style = "Normal"

$wsSummary.Range("J6").NumberFormat = "@"
$wsSummary.Range("J6").Value = "16.3"
$wsSummary.Range("J6").Style = "Normal"

$wsSummary.Range("L6").NumberFormat = "@"
$wsSummary.Range("L6").Value = "54.7"
$wsSummary.Range("L6").Style = "Normal"

$wsSummary.Range("N6").NumberFormat = "@"
$wsSummary.Range("N6").Value = "24.1"
$wsSummary.Range("N6").Style = "Normal"
